$wb = $excel.ActiveWorkbook

# --- HARDWARE_MANAGEMENT: add 4 new BMC port rows (25-28) for river compute nodes ---
$hw = $wb.Worksheets.Item("HARDWARE_MANAGEMENT")

$hw.Range("J25").Value = "cn01"
$hw.Range("K25").Value = "x3002"
$hw.Range("L25").Value = "u15"
$hw.Range("M25").Value = "bmc"
$hw.Range("O25").Value = 1
$hw.Range("P25").Value = "sw-leaf-bmc-001"
$hw.Range("Q25").Value = "x3000"
$hw.Range("R25").Value = "u37"
$hw.Range("T25").Value = 11

$hw.Range("J26").Value = "cn02"
$hw.Range("K26").Value = "x3002"
$hw.Range("L26").Value = "u16"
$hw.Range("M26").Value = "bmc"
$hw.Range("O26").Value = 1
$hw.Range("P26").Value = "sw-leaf-bmc-001"
$hw.Range("Q26").Value = "x3000"
$hw.Range("R26").Value = "u37"
$hw.Range("T26").Value = 12

$hw.Range("J27").Value = "cn03"
$hw.Range("K27").Value = "x3002"
$hw.Range("L27").Value = "u17"
$hw.Range("M27").Value = "bmc"
$hw.Range("O27").Value = 1
$hw.Range("P27").Value = "sw-leaf-bmc-001"
$hw.Range("Q27").Value = "x3000"
$hw.Range("R27").Value = "u37"
$hw.Range("T27").Value = 13

$hw.Range("J28").Value = "cn04"
$hw.Range("K28").Value = "x3002"
$hw.Range("L28").Value = "u18"
$hw.Range("M28").Value = "bmc"
$hw.Range("O28").Value = 1
$hw.Range("P28").Value = "sw-leaf-bmc-001"
$hw.Range("Q28").Value = "x3000"
$hw.Range("R28").Value = "u37"
$hw.Range("T28").Value = 14

$hw.Range("J25:T28").Select()

# --- COMPUTE_NODES: add 4 new BMC port rows (24-27) for river compute nodes ---
$cn = $wb.Worksheets.Item("COMPUTE_NODES")

$cn.Range("J24").Value = "cn01"
$cn.Range("K24").Value = "x3002"
$cn.Range("L24").Value = "u15"
$cn.Range("O24").Value = 1
$cn.Range("P24").Value = "sw-leaf-bmc-001"
$cn.Range("Q24").Value = "x3000"
$cn.Range("R24").Value = "u37"
$cn.Range("T24").Value = 24

$cn.Range("J25").Value = "cn02"
$cn.Range("K25").Value = "x3002"
$cn.Range("L25").Value = "u16"
$cn.Range("O25").Value = 1
$cn.Range("P25").Value = "sw-leaf-bmc-001"
$cn.Range("Q25").Value = "x3000"
$cn.Range("R25").Value = "u37"
$cn.Range("T25").Value = 25

$cn.Range("J26").Value = "cn03"
$cn.Range("K26").Value = "x3002"
$cn.Range("L26").Value = "u17"
$cn.Range("O26").Value = 1
$cn.Range("P26").Value = "sw-leaf-bmc-001"
$cn.Range("Q26").Value = "x3000"
$cn.Range("R26").Value = "u37"
$cn.Range("T26").Value = 26

$cn.Range("J27").Value = "cn04"
$cn.Range("K27").Value = "x3002"
$cn.Range("L27").Value = "u18"
$cn.Range("O27").Value = 1
$cn.Range("P27").Value = "sw-leaf-bmc-001"
$cn.Range("Q27").Value = "x3000"
$cn.Range("R27").Value = "u37"
$cn.Range("T27").Value = 27

$cn.Range("J24:T27").Select()
$cn.Activate()
